# Generate Report for Handoff
#
# The handoff source document's GUID changed from
#   83284c0a-188f-491c-a5ba-993e9c7321c7
# to
#   ab744e59-8a88-4337-b0c7-b9fb9cfbfb7d
# and a brand new handoff round was generated, producing new xlf package
# hashes (9e09b058d56999d7ed4227757d38e8f3ee73e258 -> 55f9867294a4dc4379d7b9715191150b774c6e91)
# and new handoff timestamps for both locales.
#
# This script updates the cell text (and the backing shared strings) as
# well as the cached hyperlink display text on every sheet so both stay in
# sync, without touching the hyperlink targets (Address/rels), which the
# original commit also left untouched.

$wb = $excel.ActiveWorkbook

$oldGuid = "83284c0a-188f-491c-a5ba-993e9c7321c7"
$newGuid = "ab744e59-8a88-4337-b0c7-b9fb9cfbfb7d"

$oldHash = "9e09b058d56999d7ed4227757d38e8f3ee73e258"
$newHash = "55f9867294a4dc4379d7b9715191150b774c6e91"

$newMdName = "$newGuid.md"

$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$newZhHandoffTime = "2016-03-10 21:07:49"
$newDeHandoffTime = "2016-03-10 21:07:55"

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Update the cell text (drives the shared-string table) -----------------

$ws1.Range("A2").Value = $newMdName

$ws2.Range("A2").Value = $newMdName
$ws2.Range("C2").Value = $newZhXlf
$ws2.Range("D2").Value = $newZhHandoffTime

$ws3.Range("A2").Value = $newMdName
$ws3.Range("C2").Value = $newDeXlf
$ws3.Range("D2").Value = $newDeHandoffTime

# --- Update the cached hyperlink display text in place ---------------------
# Iterating with foreach binds each Hyperlink object to its own entry so the
# update happens in place (same r:id / target, no duplicate link records).

foreach ($h in $ws1.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    }
}

foreach ($h in $ws2.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = $newZhXlf
    }
}

foreach ($h in $ws3.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = $newDeXlf
    }
}
